$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (progress for "App" / week 4) values
$ws.Range("E2").Value = 100
$ws.Range("E3").Value = 100
$ws.Range("E4").Value = 88
$ws.Range("E5").Value = 100
$ws.Range("E6").Value = 100
$ws.Range("E7").Value = 100
$ws.Range("E8").Value = 100
$ws.Range("E9").Value = 88
$ws.Range("E10").Value = 91
$ws.Range("E11").Value = 90
$ws.Range("E12").Value = 70
$ws.Range("E13").Value = 60
$ws.Range("E14").Value = 100

# Updates to column D
$ws.Range("D11").Value = 90
$ws.Range("D12").Value = 70
$ws.Range("D13").Value = 60

# Update to column C
$ws.Range("C13").Value = 60

# Update the active cell / selection in the sheet view
$ws.Range("E1").Select()
